$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "newSheet"

$ws.Range("A4").Value = "PBSSD"
$ws.Range("B4").Value = "MTBF"

function Set-FullBox($rng) {
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(7).Color = 0
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(8).Color = 0
  $rng.Borders.Item(9).LineStyle = 1
  $rng.Borders.Item(9).Color = 0
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(10).Color = 0
}

$a4 = $ws.Range("A4")
Set-FullBox $a4

$b4 = $ws.Range("B4")
$b4.Borders.Item(8).LineStyle = 1
$b4.Borders.Item(8).Color = 0
$b4.Borders.Item(9).LineStyle = 1
$b4.Borders.Item(9).Color = 0
$b4.Borders.Item(10).LineStyle = 1
$b4.Borders.Item(10).Color = 0

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $last)
